# Apply updated crypto price/volume data per commit "Updated cryptos list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their original plain-text representation
# (avoids Excel auto-converting numeric-looking strings like "1.00" or "0.0672"
# into real numbers, which would drop formatting such as trailing zeros).
$ws.Range("D2").Value = "67.334.98"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "3.325.32"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "186.53"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "578.25"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.407"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "3.893.66"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "67.538.54"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "3.324.71"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "444.46"
$ws.Range("E18").Value = "  +7.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.68"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.58"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.73"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.20"
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.463.03"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.516"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.97"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.94"
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").Value = "  +5.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.68"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "2.785.77"
$ws.Range("E40").Value = "  +5.62%  "
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.48"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.24"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0672"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.80"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "326.72"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.991"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.18"
$ws.Range("E51").Value = "  +1.77%  "
